# "fin du xhi deux"
# Corrects the "sexe" (column C) values for several subjects from "f" to "h"
# on rows 25, 26, 27, 35, 43, 44 and 45, and leaves the selection on the
# last cell that was touched/reviewed (E33), matching where the author
# ended up scrolling/looking in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToFix = @(25, 26, 27, 35, 43, 44, 45)
foreach ($r in $rowsToFix) {
    $ws.Range("C$r").Value = "h"
}

# Reflect the final cursor/selection position seen in the saved workbook.
$ws.Range("E33").Select()
